# Recalibration update: refresh PnL factor-model calibration outputs
# (risk driver type PnL) across the AR, SETAR, GARCH, TARCH and AR_TARCH
# worksheets, per "Settings back to qfin22 presentation - recalibration".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.0001937795527548159
$ws.Range("B3").Value = 0.7960856398156899
$ws.Range("B4").Value = 0.001479385803919736
$ws.Range("B5").Value = "[0.9999999999999999, 0.33946158343711, 0.34125169465527905, 0.3188024158459753, 0.30377310294201115, 0.501238080154443, 0.32595732762256935, 0.33342891480635306, 0.3196251240418694, 0.30049901821935254, 0.2946029839065674, 0.28121225603548694, 0.2909792475500102, 0.31329501186246744, 0.32948512612660086, 0.3026650722816156, 0.2772984323524291, 0.29193030018116306, 0.28908737459548817, 0.28623648591838696]"

# ---------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.009358913782744869
$ws.Range("B5").Value = 0.7698223186380851
$ws.Range("B6").Value = 0.0009379016832826149
$ws.Range("B7").Value = 0.009780399288728352
$ws.Range("B8").Value = 0.7607727225675658
$ws.Range("B9").Value = 0.001023018040085841
$ws.Range("B10").Value = "[0.9999999999999999, 0.36445934353256837, 0.34174219069493267, 0.34126403311521153, 0.30389410785702076, 0.3732957348543778, 0.3194669107101271, 0.3114748741393112, 0.29904154720229803, 0.27547340248766733, 0.30015314165955126, 0.26792893101864956, 0.26482126499654324, 0.2635847307644695, 0.2885620742518789, 0.31420336433637375, 0.26485914739532673, 0.26640756189680664, 0.27879395384771616, 0.26348342597892405]"

# ---------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.0003279060636946329
$ws.Range("B3").Value = [double]"3.295211686320229E-05"
$ws.Range("B4").Value = 0.2000062398467103
$ws.Range("B5").Value = 0.7800028897075
$ws.Range("B6").Value = "[1.0, -0.08331080408483663, -0.05528697953383657, -0.09718775020770687, -0.09252378413800325, 0.3013669140283891, -0.03309454254153107, 0.05686365641186162, -0.0032054331467312237, -0.006653781691703549, -0.02632601402050589, -0.025754428478049724, -0.0008811668881979709, 0.033030436795600734, 0.08138672992058991, 0.02853514013493564, -0.0017776299991643057, -0.010573201985815473, -0.01514513509068795, 0.025210419763523204]"

# ---------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0005020158321439698
$ws.Range("B3").Value = [double]"3.295198358091907E-05"
$ws.Range("B4").Value = 0.1000002843557148
$ws.Range("B5").Value = 0.8299997465973533
$ws.Range("B6").Value = "[1.0, -0.05073861929783802, -0.02965419603942369, -0.07551127133353605, -0.07020024626513445, 0.3124014038565995, -0.020414943355307802, 0.06732504184971372, 0.0073206729473030375, -0.00025653277107786645, -0.01829765302055078, -0.021379580273064587, 0.005446280220337026, 0.038991993914304475, 0.08335169487976186, 0.034807924912723776, 0.0001596533967078244, -0.005335192729075716, -0.01075128959129605, 0.030967311475503696]"
$ws.Range("B7").Value = 0.09999981829131943

# ---------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = [double]"7.139045213874651E-06"
$ws.Range("B3").Value = [double]"2.9858546605661E-05"
$ws.Range("B4").Value = 0.1999963662190925
$ws.Range("B5").Value = 0.7749858006593405
$ws.Range("B6").Value = "[0.9999999999999998, -0.06810290893194132, -0.03310670119976743, -0.0679405337271679, -0.049842397424505816, 0.22968097745468785, -0.025130993354099157, 0.036645467501377704, -0.006245278606965087, 0.0014279777449291054, -0.012995261281504134, -0.014942720531516046, -0.0035473078309468916, 0.02193623646241282, 0.07325091140556408, 0.030698685884510996, -0.0011529299312756168, 0.0038486431419186245, 0.00826934799085138, 0.02151196293244051]"
$ws.Range("B7").Value = 0.009997765824648007
$ws.Range("B9").Value = 0.7896263652286247

$wb.Save()
